# Updates cryptos price/volume figures to the latest scrape, and re-orders
# the BabyDogeCoin/Quant and EnergySwap/Algorand rows to match new rankings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    if ($val -match '^-?\d+(\.\d+)?$') {
        # Looks like a plain number (e.g. "219.90") - Excel would silently
        # coerce it to a numeric value and drop formatting (trailing zeros).
        # Force text via a leading quote-prefix, then strip the resulting
        # quotePrefix style so the cell keeps the sheet default style.
        $cell.Value = "'" + $val
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $val
    }
}

Set-TextCell 2 4 '28.300.75'
Set-TextCell 2 5 '  +4.00%  '
Set-TextCell 3 4 '1.733.47'
Set-TextCell 3 5 '  +2.99%  '
Set-TextCell 4 5 '  +0.10%  '
Set-TextCell 5 4 '219.90'
Set-TextCell 5 5 '  +1.98%  '
Set-TextCell 6 5 '  +0.90%  '
Set-TextCell 7 5 '  +0.01%  '
Set-TextCell 8 4 '24.13'
Set-TextCell 8 5 '  +11.53%  '
Set-TextCell 9 4 '0.269'
Set-TextCell 9 5 '  +4.78%  '
Set-TextCell 10 4 '0.0639'
Set-TextCell 10 5 '  +2.36%  '
Set-TextCell 11 5 '  +0.75%  '
Set-TextCell 12 4 '1.978.47'
Set-TextCell 12 5 '  +3.05%  '
Set-TextCell 13 4 '1.733.40'
Set-TextCell 13 5 '  +3.06%  '
Set-TextCell 14 4 '4.28'
Set-TextCell 14 5 '  +2.95%  '
Set-TextCell 15 5 '  +4.08%  '
Set-TextCell 16 4 '67.73'
Set-TextCell 16 5 '  +1.97%  '
Set-TextCell 17 4 '28.279.11'
Set-TextCell 17 5 '  +4.02%  '
Set-TextCell 18 4 '242.32'
Set-TextCell 18 5 '  +1.48%  '
Set-TextCell 19 5 '  +1.77%  '
Set-TextCell 20 4 '7.98'
Set-TextCell 20 5 '  -1.23%  '
Set-TextCell 21 5 '  +0.07%  '
Set-TextCell 22 5 '  +2.31%  '
Set-TextCell 23 4 '9.70'
Set-TextCell 23 5 '  +2.18%  '
Set-TextCell 24 5 '  +0.21%  '
Set-TextCell 25 4 '150.03'
Set-TextCell 25 5 '  +1.24%  '
Set-TextCell 26 4 '7.54'
Set-TextCell 26 5 '  +3.71%  '
Set-TextCell 27 4 '16.69'
Set-TextCell 27 5 '  +2.08%  '
Set-TextCell 28 5 '  +1.13%  '
Set-TextCell 29 5 '  +0.09%  '
Set-TextCell 30 4 '0.0515'
Set-TextCell 30 5 '  +3.11%  '
Set-TextCell 31 5 '  +2.60%  '
Set-TextCell 32 4 '3.45'
Set-TextCell 32 5 '  +2.08%  '
Set-TextCell 33 4 '1.505.81'
Set-TextCell 33 5 '  -4.24%  '
Set-TextCell 34 5 '  +1.36%  '
Set-TextCell 35 5 '  -1.99%  '
Set-TextCell 36 4 '0.967'
Set-TextCell 36 5 '  +2.84%  '
Set-TextCell 37 5 '  +0.64%  '
Set-TextCell 38 4 '2.41'
Set-TextCell 38 5 '  +0.61%  '
Set-TextCell 39 5 '  +1.84%  '
Set-TextCell 40 4 '1.08'
Set-TextCell 40 5 '  +1.51%  '
Set-TextCell 41 5 '  +2.10%  '
Set-TextCell 42 5 '  +1.57%  '
Set-TextCell 43 5 '  +0.01%  '
Set-TextCell 44 5 '  +2.11%  '
Set-TextCell 45 4 '1.881.91'
Set-TextCell 45 5 '  +2.81%  '
Set-TextCell 46 4 '0.803'
Set-TextCell 46 5 '  +1.82%  '
Set-TextCell 47 5 '  +8.96%  '

# Rows 48-51 were re-ranked: BabyDogeCoin/Quant swap places (rows 48-49),
# and EnergySwap/Algorand swap places (rows 50-51). Each swap also carries
# updated price + volume figures, so the full B:E block is rewritten.

Set-TextCell 48 2 'Quant'
Set-TextCell 48 3 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 48 4 '91.03'
Set-TextCell 48 5 '  +0.08%  '

Set-TextCell 49 2 'BabyDogeCoin'
Set-TextCell 49 3 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell 49 4 '0.0₆0114'
Set-TextCell 49 5 '  +7.37%  '

Set-TextCell 50 2 'Algorand'
Set-TextCell 50 3 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 50 4 '0.105'
Set-TextCell 50 5 '  +0.71%  '

Set-TextCell 51 2 'EnergySwap'
Set-TextCell 51 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 51 4 '8.21'
Set-TextCell 51 5 '  +0.50%  '
